$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.122973680496216
$ws.Range("B1").Value = 2.752149820327759
$ws.Range("C1").Value = 1.863809466362
$ws.Range("D1").Value = 1.524407982826233
$ws.Range("E1").Value = 1.425798296928406
